# feat: filter products by tag/category
# Adds a new "category_id" column (C) to the products sheet.
# Rows 2-60 (seller_id = 1 or 3) get category_id = 2.
# Rows 61-102 (seller_id = 24, second block) get category_id = 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Cells.Item(1, 3).Value = "category_id"

# Body: rows 2..60 -> 2, rows 61..102 -> 5
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 3).Value = 2
}
for ($r = 61; $r -le 102; $r++) {
    $ws.Cells.Item($r, 3).Value = 5
}

# Keep the selection/view roughly where the author left it.
$ws.Range("G104").Select()
